# Update countries & provincias Spain
# - Swap country order for Egipto/Luxemburgo (Egipto overtook Luxemburgo in case count)
# - Swap country order for Sri Lanka/Mauricio (Sri Lanka overtook Mauricio in case count)
# - Update "Datos actualizados" timestamp
# - Refresh case counts for several countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country label swaps (row 53/54: Luxemburgo <-> Egipto) ---
$ws.Range("A53").Value = "Egipto"
$ws.Range("A54").Value = "Luxemburgo"

# --- Country label swaps (row 113/114: Mauricio <-> Sri Lanka) ---
$ws.Range("A113").Value = "Sri Lanka"
$ws.Range("A114").Value = "Mauricio"

# --- Last updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 20:22"

# --- Row 7: Francia ---
$ws.Range("B7").Value = 159877
$ws.Range("C7").Value = 1827
$ws.Range("D7").Value = 40657
$ws.Range("E7").Value = 97880
$ws.Range("F7").Value = 5218
$ws.Range("G7").Value = 544
$ws.Range("H7").Value = 21340

# --- Row 25: Israel ---
$ws.Range("B25").Value = 14498
$ws.Range("C25").Value = 556
$ws.Range("D25").Value = 5215
$ws.Range("E25").Value = 9094
$ws.Range("F25").Value = 141
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 189

# --- Row 40: Noruega ---
$ws.Range("E40").Value = 7056
$ws.Range("G40").Value = 5
$ws.Range("H40").Value = 187

# --- Row 53: Egipto (new data, after swap) ---
$ws.Range("B53").Value = 3659
$ws.Range("C53").Value = 169
$ws.Range("D53").Value = 935
$ws.Range("E53").Value = 2448
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 12
$ws.Range("H53").Value = 276

# --- Row 54: Luxemburgo (old data, after swap) ---
$ws.Range("B54").Value = 3654
$ws.Range("C54").Value = 36
$ws.Range("D54").Value = 711
$ws.Range("E54").Value = 2863
$ws.Range("F54").Value = 32
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 80

# --- Row 56: Marruecos ---
$ws.Range("B56").Value = 3446
$ws.Range("C56").Value = 237
$ws.Range("D56").Value = 417
$ws.Range("E56").Value = 2880

# --- Row 65: Barein ---
$ws.Range("B65").Value = 2027
$ws.Range("C65").Value = 54
$ws.Range("E65").Value = 994

# --- Row 108: Jordania ---
$ws.Range("D108").Value = 315
$ws.Range("E108").Value = 113

# --- Row 113: Sri Lanka (new data, after swap) ---
$ws.Range("B113").Value = 330
$ws.Range("C113").Value = 20
$ws.Range("D113").Value = 105
$ws.Range("E113").Value = 218
$ws.Range("F113").Value = 2
$ws.Range("H113").Value = 7

# --- Row 114: Mauricio (old data, after swap) ---
$ws.Range("B114").Value = 329
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = 261
$ws.Range("E114").Value = 59
$ws.Range("F114").Value = 3
$ws.Range("H114").Value = 9

# --- Row 149: Guinea Ecuatorial ---
$ws.Range("B149").Value = 84
$ws.Range("C149").Value = 1
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 1

# --- Row 154: San Martin (Parte Holandesa) ---
$ws.Range("B154").Value = 71
$ws.Range("C154").Value = 3
$ws.Range("D154").Value = 22
$ws.Range("E154").Value = 38
$ws.Range("F154").Value = 2
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 11
